# Hortaliza, Femacal de La Calera - Ají
# A new price record is inserted as row 610 (pushing the existing rows
# 610..698 down to 611..699), adding one more weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 610, shifting rows 610:698 down to 611:699.
$ws.Rows.Item(610).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(610, 1).Value  = 3
$ws.Cells.Item(610, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(610, 3).Value  = "Coquimbo"
$ws.Cells.Item(610, 4).Value  = 45077
$ws.Cells.Item(610, 5).Value  = 5
$ws.Cells.Item(610, 6).Value  = 100112021
$ws.Cells.Item(610, 7).Value  = "Ají"
$ws.Cells.Item(610, 8).Value  = "Inferno"
$ws.Cells.Item(610, 9).Value  = "Primera"
$ws.Cells.Item(610, 10).Value = 85
$ws.Cells.Item(610, 11).Value = 13000
$ws.Cells.Item(610, 12).Value = 14000
$ws.Cells.Item(610, 13).Value = 13471
$ws.Cells.Item(610, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(610, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(610, 16).Value = 1347
$ws.Cells.Item(610, 17).Value = 10
$ws.Cells.Item(610, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D.
$ws.Cells.Item(610, 4).NumberFormat = $ws.Cells.Item(611, 4).NumberFormat
